# This script normalizes several shared-string values across the workbook
# that previously contained embedded line breaks, replacing the line breaks
# with single spaces so the text reads as one line. One of these values
# ("Fluvirin" / "Preservative-free") turns out to duplicate a value that
# already exists elsewhere in the workbook once the line break is removed,
# so Excel will automatically collapse/dedupe the shared string table when
# the value is rewritten.

$wb = $excel.ActiveWorkbook

$adultVaccine = $wb.Worksheets.Item("Adult Vaccine ")
$pedFlu       = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$adultFlu     = $wb.Worksheets.Item("Adult Influenza Vaccine ")

# Adult Vaccine sheet
$adultVaccine.Range("B14").Value = "Tetanus  Diphtheria Toxoids Adsorbed for Adults No Preservative"

# Pediatric Influenza Vaccine sheet
$pedFlu.Range("B3").Value  = "Fluzone Pediatric dose No Preservative"
$pedFlu.Range("B6").Value  = "Fluarix Preservative-Free"
$pedFlu.Range("B9").Value  = "FluMist No Preservative"
$pedFlu.Range("B10").Value = "Afluria No Preservative"
$pedFlu.Range("B12").Value = "Afluria No Preservative"
$pedFlu.Range("H10").Value = "Merck (CSL product)"
$pedFlu.Range("H11").Value = "Merck (CSL product)"
$pedFlu.Range("H12").Value = "Merck (CSL product)"

# Adult Influenza Vaccine sheet
$adultFlu.Range("B5").Value  = "Agriflu No Preservative"
$adultFlu.Range("B7").Value  = "Fluvirin Preservative-free"
$adultFlu.Range("B8").Value  = "Fluraix Preservative-free"
$adultFlu.Range("B10").Value = "Flumist No Preservative"
